# Scheduled runner update: refresh market-price-derived profit figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) across several
# crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 2041.3334
$ws.Range("J48").Value = 2041.3334
$ws.Range("L48").Value = 6124.0002
$ws.Range("N48").Value = -6708.0002
$ws.Range("H56").Value = 2041.3334
$ws.Range("J56").Value = 2041.3334
$ws.Range("L56").Value = 6124.0002
$ws.Range("N56").Value = -7192.0002
$ws.Range("H76").Value = 8749.5
$ws.Range("I76").Value = 8999.333000000001
$ws.Range("K76").Value = 8999.333000000001
$ws.Range("M76").Value = -8684.333000000001
$ws.Range("H79").Value = 8749.5
$ws.Range("I79").Value = 8999.333000000001
$ws.Range("K79").Value = 8999.333000000001
$ws.Range("M79").Value = -7907.333000000001
$ws.Range("H88").Value = 18942.16
$ws.Range("I88").Value = 2165.125
$ws.Range("J88").Value = 26837.234
$ws.Range("K88").Value = 2165.125
$ws.Range("L88").Value = 26837.234
$ws.Range("M88").Value = -1759.125
$ws.Range("N88").Value = -27649.234
$ws.Range("H91").Value = 18942.16
$ws.Range("I91").Value = 2165.125
$ws.Range("J91").Value = 26837.234
$ws.Range("K91").Value = 2165.125
$ws.Range("L91").Value = 26837.234
$ws.Range("M91").Value = -761.125
$ws.Range("N91").Value = -29645.234
$ws.Range("H98").Value = 52538.145
$ws.Range("I98").Value = 31261.268
$ws.Range("J98").Value = 180199.4
$ws.Range("K98").Value = 31261.268
$ws.Range("L98").Value = 180199.4
$ws.Range("M98").Value = -29763.268
$ws.Range("N98").Value = -183195.4
$ws.Range("H101").Value = 3395.8
$ws.Range("J101").Value = 3484.6667
$ws.Range("L101").Value = 10454.0001
$ws.Range("N101").Value = -13698.0001
$ws.Range("H104").Value = 154.875
$ws.Range("I104").Value = 154.875
$ws.Range("K104").Value = 464.625
$ws.Range("M104").Value = 1282.375
$ws.Range("H122").Value = 52538.145
$ws.Range("I122").Value = 31261.268
$ws.Range("J122").Value = 180199.4
$ws.Range("K122").Value = 93783.804
$ws.Range("L122").Value = 540598.2
$ws.Range("M122").Value = -91333.804
$ws.Range("N122").Value = -545498.2
$ws.Range("H138").Value = 3003.1177
$ws.Range("J138").Value = 4463.0713
$ws.Range("L138").Value = 13389.2139
$ws.Range("N138").Value = -23669.2139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 78850
$ws.Range("J103").Value = 78850
$ws.Range("L103").Value = 78850
$ws.Range("N103").Value = -81194
$ws.Range("H110").Value = 29095.455
$ws.Range("I110").Value = 32281.207
$ws.Range("J110").Value = 5998.75
$ws.Range("K110").Value = 32281.207
$ws.Range("L110").Value = 5998.75
$ws.Range("M110").Value = -30236.207
$ws.Range("N110").Value = -10088.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2138.35
$ws.Range("I105").Value = 1775
$ws.Range("K105").Value = 1775
$ws.Range("M105").Value = -28

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 16931.666
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 16931.666
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 16931.666
$ws.Range("N37").Value = -17145.666
$ws.Range("M37").ClearContents()
$ws.Range("H59").Value = 46537.125
$ws.Range("J59").Value = 50328.145
$ws.Range("L59").Value = 50328.145
$ws.Range("N59").Value = -52618.145
$ws.Range("H62").Value = 78741.57000000001
$ws.Range("I62").Value = 205061
$ws.Range("J62").Value = 8564.111000000001
$ws.Range("K62").Value = 205061
$ws.Range("L62").Value = 8564.111000000001
$ws.Range("M62").Value = -204437
$ws.Range("N62").Value = -9812.111000000001
$ws.Range("H65").Value = 78741.57000000001
$ws.Range("I65").Value = 205061
$ws.Range("J65").Value = 8564.111000000001
$ws.Range("K65").Value = 1025305
$ws.Range("L65").Value = 42820.55500000001
$ws.Range("M65").Value = -1022185
$ws.Range("N65").Value = -49060.55500000001
$ws.Range("H105").Value = 1189.6666
$ws.Range("I105").Value = 1227.6
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 1227.6
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 519.4000000000001
$ws.Range("N105").Value = -4494
$ws.Range("H134").Value = 4427.25
$ws.Range("I134").Value = 4140.4736
$ws.Range("K134").Value = 12421.4208
$ws.Range("M134").Value = -9886.4208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 628.5
$ws.Range("I8").Value = 628.5
$ws.Range("K8").Value = 1885.5
$ws.Range("M8").Value = -1746.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 344000
$ws.Range("J95").Value = 344000
$ws.Range("L95").Value = 344000
$ws.Range("N95").Value = -349492
$ws.Range("H122").Value = 3078.9697
$ws.Range("I122").Value = 2338.3809
$ws.Range("K122").Value = 7015.1427
$ws.Range("M122").Value = -4565.1427

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2645.6667
$ws.Range("I16").Value = 1424.8
$ws.Range("K16").Value = 1424.8
$ws.Range("M16").Value = -1254.8
$ws.Range("H40").Value = 35288.844
$ws.Range("I40").Value = 39327.066
$ws.Range("J40").Value = 11059.5
$ws.Range("K40").Value = 39327.066
$ws.Range("L40").Value = 11059.5
$ws.Range("M40").Value = -39191.066
$ws.Range("N40").Value = -11331.5
$ws.Range("H82").Value = 3500.4285
$ws.Range("J82").Value = 5003
$ws.Range("L82").Value = 5003
$ws.Range("N82").Value = -5725
$ws.Range("H85").Value = 3500.4285
$ws.Range("J85").Value = 5003
$ws.Range("L85").Value = 5003
$ws.Range("N85").Value = -7499
$ws.Range("H103").Value = 36031.383
$ws.Range("J103").Value = 36867.332
$ws.Range("L103").Value = 36867.332
$ws.Range("N103").Value = -39211.332
$ws.Range("H136").Value = 4386.6
$ws.Range("I136").Value = 3762.889
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 11288.667
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -8738.667000000001
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 12515248
$ws.Range("I4").Value = 16678666
$ws.Range("J4").Value = 24994
$ws.Range("K4").Value = 16678666
$ws.Range("L4").Value = 24994
$ws.Range("M4").Value = -16678553
$ws.Range("N4").Value = -25220
$ws.Range("H62").Value = 116960.89
$ws.Range("J62").Value = 7387.778
$ws.Range("L62").Value = 7387.778
$ws.Range("N62").Value = -8635.778
$ws.Range("H65").Value = 116960.89
$ws.Range("J65").Value = 7387.778
$ws.Range("L65").Value = 36938.89
$ws.Range("N65").Value = -43178.89
$ws.Range("H69").Value = 50663.332
$ws.Range("J69").Value = 50663.332
$ws.Range("L69").Value = 50663.332
$ws.Range("N69").Value = -52161.332
$ws.Range("H72").Value = 50663.332
$ws.Range("J72").Value = 50663.332
$ws.Range("L72").Value = 151989.996
$ws.Range("N72").Value = -159477.996
$ws.Range("H100").Value = 1331.1305
$ws.Range("I100").Value = 1246.5
$ws.Range("K100").Value = 2493
$ws.Range("M100").Value = -1952
$ws.Range("H101").Value = 12583.333
$ws.Range("J101").Value = 12583.333
$ws.Range("L101").Value = 12583.333
$ws.Range("N101").Value = -19073.333
$ws.Range("H136").Value = 3943.6538
$ws.Range("I136").Value = 3441.739
$ws.Range("J136").Value = 7791.6665
$ws.Range("K136").Value = 10325.217
$ws.Range("L136").Value = 23374.9995
$ws.Range("M136").Value = -7775.217000000001
$ws.Range("N136").Value = -28474.9995
